# "Fruta / hortaliza, semanal" — weekly refresh of the Acelga (Vega Modelo
# de Temuco) price series: a new week's observation is inserted at row 114
# (pushing the existing rows 114:219 down to 115:220) and populated with
# the latest reading; everything else keeps its prior values/order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 114, shifting rows 114-219 down
# to 115-220 (dimension grows from R219 to R220).
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the new week's data. The fixed
# descriptive columns (market/region/category/etc.) repeat the same values
# as the rest of the series.
$ws.Cells.Item(114, 1).Value = 10
$ws.Cells.Item(114, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(114, 3).Value = "La Araucanía"
$ws.Cells.Item(114, 4).Value = 44512
$ws.Cells.Item(114, 5).Value = 9
$ws.Cells.Item(114, 6).Value = 100112009
$ws.Cells.Item(114, 7).Value = "Acelga"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 30
$ws.Cells.Item(114, 11).Value = 9000
$ws.Cells.Item(114, 12).Value = 9000
$ws.Cells.Item(114, 13).Value = 9000
$ws.Cells.Item(114, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(114, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(114, 16).Value = 750
$ws.Cells.Item(114, 17).Value = 12
$ws.Cells.Item(114, 18).Value = "Hortaliza"
